# Apply edits to "Test Suite" worksheet (TC2 / TC3 step text corrections)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# TC2 - Step 2 (row 21): action text changes from "seleciona..." to "preenche os campos..."
$ws.Range("B21").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"

# TC2 - Step 2 (row 21): expected result changes from "TJSeg..." to "nome de usuario e/ou senha estao incorretos"
$ws.Range("D21").Value = "SYSTEM alerta que o nome de usuario e/ou senha estao incorretos"

# TC3 - Step 2 (row 31): expected result changes from "nome de usuario e/ou senha estao incorretos" to "TJSeg..."
$ws.Range("D31").Value = "SYSTEM alerta que o TJSeg (sistema que fornece as permissoes de acesso e escrita) esta fora do ar"
